$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 7).Value = "Aansluiting"  # G1
$ws.Cells.Item(1, 8).Value = "Vooropleiding"  # H1
$ws.Cells.Item(2, 7).Value = "Direct"  # G2
$ws.Cells.Item(2, 8).Value = "HAVO"  # H2
$ws.Cells.Item(3, 7).Value = "Direct"  # G3
$ws.Cells.Item(3, 8).Value = "HAVO"  # H3
$ws.Cells.Item(4, 1).Value = "Aansluiting"  # A4
$ws.Cells.Item(4, 2).Value = "Direct"  # B4
$ws.Cells.Item(4, 4).Value = 756  # D4
$ws.Cells.Item(4, 5).Value = 0.469  # E4
$ws.Cells.Item(4, 7).Value = "Direct"  # G4
$ws.Cells.Item(4, 8).Value = "HAVO"  # H4
$ws.Cells.Item(4, 9).Value = "MedV"  # I4
$ws.Cells.Item(4, 10).Value = "Nee"  # J4
$ws.Cells.Item(4, 18).Value = 19  # R4
$ws.Cells.Item(4, 19).Value = 137.5  # S4
$ws.Cells.Item(4, 20).Value = 34.2  # T4
$ws.Cells.Item(4, 27).Value = 0.1  # AA4
$ws.Cells.Item(4, 30).Value = 2016  # AD4
$ws.Cells.Item(5, 1).Value = "Aansluiting"  # A5
$ws.Cells.Item(5, 2).Value = "Tussenjaar"  # B5
$ws.Cells.Item(5, 4).Value = 176  # D5
$ws.Cells.Item(5, 5).Value = 0.109  # E5
$ws.Cells.Item(5, 6).Value = "V"  # F5
$ws.Cells.Item(5, 7).Value = "Tussenjaar"  # G5
$ws.Cells.Item(5, 8).Value = "HAVO"  # H5
$ws.Cells.Item(5, 9).Value = "MedV"  # I5
$ws.Cells.Item(5, 18).Value = 20  # R5
$ws.Cells.Item(5, 19).Value = 134.5  # S5
$ws.Cells.Item(5, 20).Value = 34.7  # T5
$ws.Cells.Item(5, 27).Value = 0.1  # AA5
$ws.Cells.Item(5, 30).Value = 2017  # AD5
$ws.Cells.Item(6, 1).Value = "Aansluiting"  # A6
$ws.Cells.Item(6, 2).Value = "Switch intern"  # B6
$ws.Cells.Item(6, 4).Value = 202  # D6
$ws.Cells.Item(6, 5).Value = 0.125  # E6
$ws.Cells.Item(6, 7).Value = "Switch intern"  # G6
$ws.Cells.Item(6, 8).Value = "HAVO"  # H6
$ws.Cells.Item(6, 15).Value = "Ja"  # O6
$ws.Cells.Item(6, 18).Value = 21  # R6
$ws.Cells.Item(6, 19).Value = 128.5  # S6
$ws.Cells.Item(6, 20).Value = 28.1  # T6
$ws.Cells.Item(6, 27).Value = 0  # AA6
$ws.Cells.Item(6, 30).Value = 2018  # AD6
$ws.Cells.Item(7, 1).Value = "Aansluiting"  # A7
$ws.Cells.Item(7, 2).Value = "Switch extern"  # B7
$ws.Cells.Item(7, 4).Value = 445  # D7
$ws.Cells.Item(7, 5).Value = 0.276  # E7
$ws.Cells.Item(7, 6).Value = "M"  # F7
$ws.Cells.Item(7, 7).Value = "Switch extern"  # G7
$ws.Cells.Item(7, 8).Value = "HAVO"  # H7
$ws.Cells.Item(7, 9).Value = "EM"  # I7
$ws.Cells.Item(7, 19).Value = 112  # S7
$ws.Cells.Item(7, 20).Value = 36.1  # T7
$ws.Cells.Item(7, 27).Value = 0  # AA7
$ws.Cells.Item(7, 28).Value = 0  # AB7
$ws.Cells.Item(8, 1).Value = "Aansluiting"  # A8
$ws.Cells.Item(8, 2).Value = "2e Studie"  # B8
$ws.Cells.Item(8, 4).Value = 15  # D8
$ws.Cells.Item(8, 5).Value = 0.009  # E8
$ws.Cells.Item(8, 7).Value = "2e Studie"  # G8
$ws.Cells.Item(8, 8).Value = "HAVO"  # H8
$ws.Cells.Item(8, 9).Value = "EM"  # I8
$ws.Cells.Item(8, 18).Value = 22  # R8
$ws.Cells.Item(8, 19).Value = 60  # S8
$ws.Cells.Item(8, 20).Value = 29.8  # T8
$ws.Cells.Item(8, 30).Value = 2019  # AD8
$ws.Cells.Item(9, 1).Value = "Aansluiting"  # A9
$ws.Cells.Item(9, 2).Value = "Na CD"  # B9
$ws.Cells.Item(9, 4).Value = 19  # D9
$ws.Cells.Item(9, 5).Value = 0.012  # E9
$ws.Cells.Item(9, 7).Value = "Na CD"  # G9
$ws.Cells.Item(9, 8).Value = "CD"  # H9
$ws.Cells.Item(9, 9).Value = "EM&CM"  # I9
$ws.Cells.Item(9, 11).Value = "Ja"  # K9
$ws.Cells.Item(9, 12).Value = "Ja"  # L9
$ws.Cells.Item(9, 13).Value = "Ja"  # M9
$ws.Cells.Item(9, 14).Value = "Ja"  # N9
$ws.Cells.Item(9, 15).Value = "Ja"  # O9
$ws.Cells.Item(9, 18).Value = 22  # R9
$ws.Cells.Item(9, 19).Value = 129  # S9
$ws.Cells.Item(9, 20).Value = 30.5  # T9
$ws.Cells.Item(9, 27).Value = 0  # AA9
$ws.Cells.Item(9, 30).Value = 2012  # AD9
$ws.Cells.Item(10, 1).Value = "Vooropleiding"  # A10
$ws.Cells.Item(10, 2).Value = "MBO"  # B10
$ws.Cells.Item(10, 4).Value = 522  # D10
$ws.Cells.Item(10, 5).Value = 0.324  # E10
$ws.Cells.Item(10, 7).Value = "Direct"  # G10
$ws.Cells.Item(10, 8).Value = "MBO"  # H10
$ws.Cells.Item(10, 9).Value = "MedV"  # I10
$ws.Cells.Item(10, 18).Value = 21  # R10
$ws.Cells.Item(10, 19).Value = 127  # S10
$ws.Cells.Item(10, 20).Value = 32.7  # T10
$ws.Cells.Item(10, 30).Value = 2017  # AD10
$ws.Cells.Item(11, 1).Value = "Vooropleiding"  # A11
$ws.Cells.Item(11, 2).Value = "HAVO"  # B11
$ws.Cells.Item(11, 4).Value = 860  # D11
$ws.Cells.Item(11, 5).Value = 0.533  # E11
$ws.Cells.Item(11, 7).Value = "Direct"  # G11
$ws.Cells.Item(11, 8).Value = "HAVO"  # H11
$ws.Cells.Item(11, 9).Value = "EM&CM"  # I11
$ws.Cells.Item(11, 11).Value = "Nee"  # K11
$ws.Cells.Item(11, 12).Value = "Nee"  # L11
$ws.Cells.Item(11, 13).Value = "Nee"  # M11
$ws.Cells.Item(11, 14).Value = "Nee"  # N11
$ws.Cells.Item(11, 15).Value = "Nee"  # O11
$ws.Cells.Item(11, 19).Value = 132  # S11
$ws.Cells.Item(11, 20).Value = 35.9  # T11
$ws.Cells.Item(12, 1).Value = "Vooropleiding"  # A12
$ws.Cells.Item(12, 2).Value = "VWO"  # B12
$ws.Cells.Item(12, 4).Value = 58  # D12
$ws.Cells.Item(12, 5).Value = 0.036  # E12
$ws.Cells.Item(12, 7).Value = "Switch extern"  # G12
$ws.Cells.Item(12, 8).Value = "VWO"  # H12
$ws.Cells.Item(12, 9).Value = "EM"  # I12
$ws.Cells.Item(12, 11).Value = "Nee"  # K12
$ws.Cells.Item(12, 12).Value = "Nee"  # L12
$ws.Cells.Item(12, 13).Value = "Nee"  # M12
$ws.Cells.Item(12, 14).Value = "Nee"  # N12
$ws.Cells.Item(12, 15).Value = "Nee"  # O12
$ws.Cells.Item(12, 18).Value = 19  # R12
$ws.Cells.Item(12, 19).Value = 132  # S12
$ws.Cells.Item(12, 20).Value = 30.9  # T12
$ws.Cells.Item(12, 27).Value = 0.1  # AA12
$ws.Cells.Item(12, 30).Value = 2016  # AD12
$ws.Cells.Item(13, 1).Value = "Vooropleiding"  # A13
$ws.Cells.Item(13, 2).Value = "BD"  # B13
$ws.Cells.Item(13, 4).Value = 92  # D13
$ws.Cells.Item(13, 5).Value = 0.057  # E13
$ws.Cells.Item(13, 7).Value = "Direct"  # G13
$ws.Cells.Item(13, 8).Value = "BD"  # H13
$ws.Cells.Item(13, 9).Value = "EM"  # I13
$ws.Cells.Item(13, 10).Value = "Onbekend"  # J13
$ws.Cells.Item(13, 11).Value = "Ja"  # K13
$ws.Cells.Item(13, 12).Value = "Ja"  # L13
$ws.Cells.Item(13, 13).Value = "Ja"  # M13
$ws.Cells.Item(13, 14).Value = "Ja"  # N13
$ws.Cells.Item(13, 15).Value = "Ja"  # O13
$ws.Cells.Item(13, 18).Value = 21  # R13
$ws.Cells.Item(13, 19).Value = 126  # S13
$ws.Cells.Item(13, 20).Value = 15.5  # T13
$ws.Cells.Item(13, 30).Value = 2017  # AD13
$ws.Cells.Item(14, 1).Value = "Vooropleiding"  # A14
$ws.Cells.Item(14, 2).Value = "CD"  # B14
$ws.Cells.Item(14, 4).Value = 30  # D14
$ws.Cells.Item(14, 5).Value = 0.019  # E14
$ws.Cells.Item(14, 7).Value = "Na CD"  # G14
$ws.Cells.Item(14, 8).Value = "CD"  # H14
$ws.Cells.Item(14, 9).Value = "EM"  # I14
$ws.Cells.Item(14, 11).Value = "Ja"  # K14
$ws.Cells.Item(14, 12).Value = "Ja"  # L14
$ws.Cells.Item(14, 13).Value = "Ja"  # M14
$ws.Cells.Item(14, 14).Value = "Ja"  # N14
$ws.Cells.Item(14, 18).Value = 22  # R14
$ws.Cells.Item(14, 19).Value = 127.5  # S14
$ws.Cells.Item(14, 20).Value = 29.6  # T14
$ws.Cells.Item(14, 30).Value = 2017.5  # AD14
$ws.Cells.Item(15, 1).Value = "Vooropleiding"  # A15
$ws.Cells.Item(15, 2).Value = "HO"  # B15
$ws.Cells.Item(15, 4).Value = 51  # D15
$ws.Cells.Item(15, 5).Value = 0.032  # E15
$ws.Cells.Item(15, 7).Value = "Switch extern"  # G15
$ws.Cells.Item(15, 8).Value = "HO"  # H15
$ws.Cells.Item(15, 9).Value = "EM"  # I15
$ws.Cells.Item(15, 11).Value = "Nee"  # K15
$ws.Cells.Item(15, 12).Value = "Nee"  # L15
$ws.Cells.Item(15, 13).Value = "Nee"  # M15
$ws.Cells.Item(15, 14).Value = "Nee"  # N15
$ws.Cells.Item(15, 15).Value = "Nee"  # O15
$ws.Cells.Item(15, 19).Value = 128  # S15
$ws.Cells.Item(15, 20).Value = 38.3  # T15
$ws.Cells.Item(15, 28).Value = 0.1  # AB15
$ws.Cells.Item(15, 30).Value = 2015  # AD15
